$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 293 (shifts existing rows 293:316 down to 294:317),
# matching the new weekly data point being prepended to the series.
$ws.Rows.Item(293).Insert()

# Populate the new row 293 with this week's data point for
# "Vega Modelo de Temuco - Zanahoria".
$ws.Range("A293").Value = 10
$ws.Range("B293").Value = "Vega Modelo de Temuco"
$ws.Range("C293").Value = "La Araucanía"
$ws.Range("D293").Value = 44746
$ws.Range("E293").Value = 9
$ws.Range("F293").Value = 100114013
$ws.Range("G293").Value = "Zanahoria"
$ws.Range("H293").Value = "Sin especificar"
$ws.Range("I293").Value = "Primera"
$ws.Range("J293").Value = 155
$ws.Range("K293").Value = 9000
$ws.Range("L293").Value = 9000
$ws.Range("M293").Value = 9000
$ws.Range("N293").Value = "$/saco 25 kilos"
$ws.Range("O293").Value = "Región de La Araucanía"
$ws.Range("P293").Value = 360
$ws.Range("Q293").Value = 25
$ws.Range("R293").Value = "Hortaliza"
